# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# 1) Bump the "Date" metadata value on the Metadata sheet.
# 2) Add a new "Mapping: Spécification métier vers l'extension ROR
#    AccessibilityLocation" column (AL) to the Elements sheet, with the
#    business mapping value "accessibiliteLieu" on the Extension.value[x]
#    row.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 (Date) -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- 2) Elements sheet: new column AL --------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Make column AL match the look & feel of the existing "Mapping: RIM Mapping"
# column (AK) - same header style, same data-row style.
$ws.Range("AK1:AK6").Copy()
$ws.Range("AL1:AL6").PasteSpecial(-4122)

$ws.Columns.Item(38).ColumnWidth = 76.77734375

$ws.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR AccessibilityLocation"
$ws.Range("AL2").Value = ""
$ws.Range("AL3").Value = ""
$ws.Range("AL4").Value = ""
$ws.Range("AL5").Value = ""
$ws.Range("AL6").Value = "accessibiliteLieu"
